# Jenkins parameterized build completed
#
# Adds a "Runmode" column (y/n/y/n) to the addCustomerTest sheet and
# makes addCustomerTest the active sheet/tab (it was test_suite before).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("addCustomerTest")

$ws.Range("E1").Value = "Runmode"
$ws.Range("E2").Value = "y"
$ws.Range("E3").Value = "n"
$ws.Range("E4").Value = "y"
$ws.Range("E5").Value = "n"

$ws.Activate()
$ws.Range("E5").Select() | Out-Null
